# Fix the typo "infroamtion" -> "information" (and drop the stray trailing
# period run) in the second paragraph of the subtitle placeholder on slide 2,
# merging the three runs into the single corrected run the diff expects.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# The second paragraph ("We discuss ... infroamtion.") starts right after the
# first paragraph's text + its paragraph-break character. Locate it robustly
# by searching for the known substring rather than hard-coding offsets.
$fullText = $tr.Text
$oldSentence = "We discuss the simple KNN (K-Nearest Neighbor) to predict the complicate problem, predicting the rating of a movie by genre and popularity infroamtion."
$newSentence = "We discuss the simple KNN (K-Nearest Neighbor) to predict the complicate problem, predicting the rating of a movie by genre and popularity information."

$startIndex = $fullText.IndexOf($oldSentence)
$target = $tr.Characters($startIndex + 1, $oldSentence.Length)
$target.Text = $newSentence
